$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 731. This shifts the old rows 731..835
# down to 732..836, matching the target diff (dimension A1:R835 -> A1:R836).
$ws.Rows(731).Insert()

# Populate the newly inserted row 731 with the new record's data.
$ws.Cells.Item(731, 1).Value = 10
$ws.Cells.Item(731, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(731, 3).Value = "La Araucanía"
$ws.Cells.Item(731, 4).Value = 45077
$ws.Cells.Item(731, 5).Value = 9
$ws.Cells.Item(731, 6).Value = 100112045
$ws.Cells.Item(731, 7).Value = "Zapallo"
$ws.Cells.Item(731, 8).Value = "Camote"
$ws.Cells.Item(731, 9).Value = "1a (guarda)"
$ws.Cells.Item(731, 10).Value = 380
$ws.Cells.Item(731, 11).Value = 500
$ws.Cells.Item(731, 12).Value = 500
$ws.Cells.Item(731, 13).Value = 500
$ws.Cells.Item(731, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(731, 15).Value = "Región del Maule"
$ws.Cells.Item(731, 16).Value = 500
$ws.Cells.Item(731, 17).Value = 1
$ws.Cells.Item(731, 18).Value = "Hortaliza"
